$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'37.045.83"
$ws.Range("D2").ClearFormats()
$ws.Range("E2").Value = "'  +1.29%  "
$ws.Range("E2").ClearFormats()
$ws.Range("D3").Value = "'2.062.02"
$ws.Range("D3").ClearFormats()
$ws.Range("E3").Value = "'  -1.93%  "
$ws.Range("E3").ClearFormats()
$ws.Range("E4").Value = "'  +0.07%  "
$ws.Range("E4").ClearFormats()
$ws.Range("D5").Value = "'249.58"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "'  -1.38%  "
$ws.Range("E5").ClearFormats()
$ws.Range("D6").Value = "'0.674"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "'  +2.54%  "
$ws.Range("E6").ClearFormats()
$ws.Range("D8").Value = "'54.56"
$ws.Range("D8").ClearFormats()
$ws.Range("E8").Value = "'  +11.72%  "
$ws.Range("E8").ClearFormats()
$ws.Range("D9").Value = "'61.10"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "'  +1.93%  "
$ws.Range("E9").ClearFormats()
$ws.Range("E10").Value = "'  +1.03%  "
$ws.Range("E10").ClearFormats()
$ws.Range("E11").Value = "'  +6.84%  "
$ws.Range("E11").ClearFormats()
$ws.Range("E12").Value = "'  +5.63%  "
$ws.Range("E12").ClearFormats()
$ws.Range("D13").Value = "'15.01"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "'  +3.03%  "
$ws.Range("E13").ClearFormats()
$ws.Range("D14").Value = "'2.362.77"
$ws.Range("D14").ClearFormats()
$ws.Range("E14").Value = "'  -1.89%  "
$ws.Range("E14").ClearFormats()
$ws.Range("E15").Value = "'  -2.29%  "
$ws.Range("E15").ClearFormats()
$ws.Range("E16").Value = "'  +4.19%  "
$ws.Range("E16").ClearFormats()
$ws.Range("D17").Value = "'2.062.74"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "'  -1.73%  "
$ws.Range("E17").ClearFormats()
$ws.Range("D18").Value = "'36.977.51"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "'  +1.11%  "
$ws.Range("E18").ClearFormats()
$ws.Range("D19").Value = "'0.0₃0933"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "'  +12.07%  "
$ws.Range("E19").ClearFormats()
$ws.Range("D20").Value = "'73.65"
$ws.Range("D20").ClearFormats()
$ws.Range("E20").Value = "'  +1.06%  "
$ws.Range("E20").ClearFormats()
$ws.Range("D21").Value = "'14.22"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "'  +6.67%  "
$ws.Range("E21").ClearFormats()
$ws.Range("D22").Value = "'5.41"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "'  +3.31%  "
$ws.Range("E22").ClearFormats()
$ws.Range("D23").Value = "'237.58"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "'  -1.23%  "
$ws.Range("E23").ClearFormats()
$ws.Range("E24").Value = "'  -0.11%  "
$ws.Range("E24").ClearFormats()
$ws.Range("E25").Value = "'  -3.89%  "
$ws.Range("E25").ClearFormats()
$ws.Range("D26").Value = "'170.03"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "'  -0.60%  "
$ws.Range("E26").ClearFormats()
$ws.Range("D27").Value = "'9.07"
$ws.Range("D27").ClearFormats()
$ws.Range("E27").Value = "'  -1.08%  "
$ws.Range("E27").ClearFormats()
$ws.Range("D28").Value = "'20.08"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "'  -5.55%  "
$ws.Range("E28").ClearFormats()
$ws.Range("D29").Value = "'2.00"
$ws.Range("D29").ClearFormats()
$ws.Range("E29").Value = "'  +0.69%  "
$ws.Range("E29").ClearFormats()
$ws.Range("E30").Value = "'  +1.98%  "
$ws.Range("E30").ClearFormats()
$ws.Range("E31").Value = "'  +2.63%  "
$ws.Range("E31").ClearFormats()
$ws.Range("E32").Value = "'  +7.23%  "
$ws.Range("E32").ClearFormats()
$ws.Range("D33").Value = "'0.0631"
$ws.Range("D33").ClearFormats()
$ws.Range("E33").Value = "'  +2.16%  "
$ws.Range("E33").ClearFormats()
$ws.Range("D34").Value = "'4.40"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "'  +7.61%  "
$ws.Range("E34").ClearFormats()
$ws.Range("D35").Value = "'0.0892"
$ws.Range("D35").ClearFormats()
$ws.Range("E35").Value = "'  -0.74%  "
$ws.Range("E35").ClearFormats()
$ws.Range("E36").Value = "'  +0.01%  "
$ws.Range("E36").ClearFormats()
$ws.Range("D37").Value = "'2.29"
$ws.Range("D37").ClearFormats()
$ws.Range("E37").Value = "'  -7.75%  "
$ws.Range("E37").ClearFormats()
$ws.Range("D38").Value = "'1.76"
$ws.Range("D38").ClearFormats()
$ws.Range("E38").Value = "'  -5.37%  "
$ws.Range("E38").ClearFormats()
$ws.Range("E39").Value = "'  +0.17%  "
$ws.Range("E39").ClearFormats()
$ws.Range("E40").Value = "'  +21.97%  "
$ws.Range("E40").ClearFormats()
$ws.Range("B41").Value = "'InjectiveProtocol"
$ws.Range("B41").ClearFormats()
$ws.Range("C41").Value = "'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj"
$ws.Range("C41").ClearFormats()
$ws.Range("D41").Value = "'17.85"
$ws.Range("D41").ClearFormats()
$ws.Range("E41").Value = "'  +10.70%  "
$ws.Range("E41").ClearFormats()
$ws.Range("B42").Value = "'VeChain"
$ws.Range("B42").ClearFormats()
$ws.Range("C42").Value = "'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("C42").ClearFormats()
$ws.Range("D42").Value = "'0.0226"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "'  +0.72%  "
$ws.Range("E42").ClearFormats()
$ws.Range("E43").Value = "'  -2.11%  "
$ws.Range("E43").ClearFormats()
$ws.Range("D44").Value = "'97.05"
$ws.Range("D44").ClearFormats()
$ws.Range("E45").Value = "'  +1.09%  "
$ws.Range("E45").ClearFormats()
$ws.Range("D46").Value = "'4.11"
$ws.Range("D46").ClearFormats()
$ws.Range("E46").Value = "'  +44.44%  "
$ws.Range("E46").ClearFormats()
$ws.Range("D47").Value = "'13.53"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "'  -52.65%  "
$ws.Range("E47").ClearFormats()
$ws.Range("D48").Value = "'2.42"
$ws.Range("D48").ClearFormats()
$ws.Range("E48").Value = "'  +8.39%  "
$ws.Range("E48").ClearFormats()
$ws.Range("D49").Value = "'1.296.73"
$ws.Range("D49").ClearFormats()
$ws.Range("E49").Value = "'  -3.16%  "
$ws.Range("E49").ClearFormats()
$ws.Range("E50").Value = "'  +2.14%  "
$ws.Range("E50").ClearFormats()
$ws.Range("D51").Value = "'4.14"
$ws.Range("D51").ClearFormats()
$ws.Range("E51").Value = "'  +7.83%  "
$ws.Range("E51").ClearFormats()
